$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K", formerly "Strike#") new values per row, regenerated per commit message.
$gValues = @{
    2  = 2
    3  = 3
    4  = 1
    5  = 0
    6  = 6
    7  = 3
    8  = 6
    9  = 5
    10 = 4
    11 = 4
    12 = 3
    13 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
